$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 11 new rows after the header row (row 1), shifting existing data down
$ws.Rows("2:12").Insert(0, 1)

# Clear the copied formatting from columns B:E (only column A should keep the date style)
$ws.Range("B2:E12").ClearFormats()

# Copy the date-number-format style used by column A into the new column-A cells
$ws.Range("A13").Copy()
$ws.Range("A2:A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the newly inserted rows with the backward-extended (1985-1995) data
$ws.Cells.Item(2, 1).Value = 31228
$ws.Cells.Item(2, 2).Value = 1985
$ws.Cells.Item(2, 3).Value = 1.058157908195101
$ws.Cells.Item(2, 4).Value = 1986
$ws.Cells.Item(2, 5).Value = -0.2746413702905404
$ws.Cells.Item(3, 1).Value = 31593
$ws.Cells.Item(3, 2).Value = 1986
$ws.Cells.Item(3, 3).Value = 0.613400232493766
$ws.Cells.Item(3, 4).Value = 1987
$ws.Cells.Item(3, 5).Value = -0.6878195673301057
$ws.Cells.Item(4, 1).Value = 31958
$ws.Cells.Item(4, 2).Value = 1987
$ws.Cells.Item(4, 3).Value = -2.809251799599333
$ws.Cells.Item(4, 4).Value = 1988
$ws.Cells.Item(4, 5).Value = -5.22296739181829
$ws.Cells.Item(5, 1).Value = 32324
$ws.Cells.Item(5, 2).Value = 1988
$ws.Cells.Item(5, 3).Value = 1.311815945521522
$ws.Cells.Item(5, 4).Value = 1989
$ws.Cells.Item(5, 5).Value = -0.2360507986455929
$ws.Cells.Item(6, 1).Value = 32689
$ws.Cells.Item(6, 2).Value = 1989
$ws.Cells.Item(6, 3).Value = 4.646141329820397
$ws.Cells.Item(6, 4).Value = 1990
$ws.Cells.Item(6, 5).Value = 4.409069002718513
$ws.Cells.Item(7, 1).Value = 33054
$ws.Cells.Item(7, 2).Value = 1990
$ws.Cells.Item(7, 3).Value = 6.234545104864941
$ws.Cells.Item(7, 4).Value = 1991
$ws.Cells.Item(7, 5).Value = 7.271550582979214
$ws.Cells.Item(8, 1).Value = 33419
$ws.Cells.Item(8, 2).Value = 1991
$ws.Cells.Item(8, 3).Value = 9.12959361095953
$ws.Cells.Item(8, 4).Value = 1992
$ws.Cells.Item(8, 5).Value = 10.06345753894098
$ws.Cells.Item(9, 1).Value = 33785
$ws.Cells.Item(9, 2).Value = 1992
$ws.Cells.Item(9, 3).Value = 3.770303686471776
$ws.Cells.Item(9, 4).Value = 1993
$ws.Cells.Item(9, 5).Value = 4.833608027690683
$ws.Cells.Item(10, 1).Value = 34150
$ws.Cells.Item(10, 2).Value = 1993
$ws.Cells.Item(10, 3).Value = -3.631379003451563
$ws.Cells.Item(10, 4).Value = 1994
$ws.Cells.Item(10, 5).Value = -4.837330717035071
$ws.Cells.Item(11, 1).Value = 34515
$ws.Cells.Item(11, 2).Value = 1994
$ws.Cells.Item(11, 3).Value = 2.525367061038386
$ws.Cells.Item(11, 4).Value = 1995
$ws.Cells.Item(11, 5).Value = 2.890876292836841
$ws.Cells.Item(12, 1).Value = 34880
$ws.Cells.Item(12, 2).Value = 1995
$ws.Cells.Item(12, 3).Value = 1.86590761958525
$ws.Cells.Item(12, 4).Value = 1996
$ws.Cells.Item(12, 5).Value = 1.408414909230937
